$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace truncated college names with their full names.
$ws.Range("B12").Value  = "Long Island University-Brooklyn Campus"
$ws.Range("B19").Value  = "North Carolina Central University"
$ws.Range("B23").Value  = "Southeastern Louisiana University"
$ws.Range("B27").Value  = "University of California, Berkeley"
$ws.Range("B32").Value  = "University of Massachusetts, Amherst"
$ws.Range("B33").Value  = "University of Missouri, Columbia"
$ws.Range("B34").Value  = "University of North Carolina Asheville"
$ws.Range("B41").Value  = "Albany State University (Georgia)"
$ws.Range("B48").Value  = "California State University, Monterey Bay"
$ws.Range("B56").Value  = "Georgia Southwestern State University"
$ws.Range("B57").Value  = "Kutztown University of Pennsylvania"
$ws.Range("B60").Value  = "Mansfield University of Pennsylvania"
$ws.Range("B61").Value  = "Metropolitan State University of Denver"
$ws.Range("B62").Value  = "Missouri Western State University"
$ws.Range("B63").Value  = "New Mexico Highlands University"
$ws.Range("B67").Value  = "Slippery Rock University of Pennsylvania"
$ws.Range("B68").Value  = "Texas A&M International University"
$ws.Range("B72").Value  = "University of Minnesota, Crookston"
$ws.Range("B74").Value  = "University of Wisconsin-Parkside"
$ws.Range("B83").Value  = "Emmanuel College (Massachusetts)"
$ws.Range("B87").Value  = "Franciscan University of Steubenville"
$ws.Range("B92").Value  = "Massachusetts College of Liberal Arts"
$ws.Range("B102").Value = "Saint Mary's University of Minnesota"
$ws.Range("B105").Value = "St. Joseph's College (Long Island)"
$ws.Range("B106").Value = "State University of New York at Geneseo"
$ws.Range("B107").Value = "State University of New York at Oswego"
$ws.Range("B108").Value = "State University of New York Maritime College"
$ws.Range("B113").Value = "University of Mary Hardin-Baylor"
$ws.Range("B121").Value = "Western Connecticut State University"
$ws.Range("B122").Value = "William Paterson University of New Jersey"

# Widen column B (names are now longer) from 11 to 21 characters.
$ws.Columns.Item(2).ColumnWidth = 20.17
